$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs", "Reln", "Lrp8", "ECs", 3, 1, 0.09865600000000001, 0.295968, 0.01009304870291239, 0.01488643315542961, 2, 0.6666666666666666, 0.7407323333333333, 2.222197, 0.125361503586392, 0.1309951108729331, 0.07307768907733334, 0.657699201696, 0.001265279761167781, 0.001950049961698009),
    @("ECs", "Reln", "Lrp8", "FAPs", 3, 1, 0.09865600000000001, 0.295968, 0.01009304870291239, 0.01488643315542961, 3, 1, 1.080588666666667, 3.241766, 0.1828787726899297, 0.1910971424199137, 0.1066065554986667, 0.9594589994880001, 0.001845804359488305, 0.002844754836827656),
    @("ECs", "Reln", "Lrp8", "M1", 3, 1, 0.09865600000000001, 0.295968, 0.01009304870291239, 0.01488643315542961, 3, 1, 1.412055666666667, 4.236167, 0.2389762314336017, 0.2497155589001607, 0.1393077638506667, 1.253769874656, 0.002411998742697805, 0.003717373975437987),
    @("ECs", "Reln", "Lrp8", "M2", 3, 1, 0.09865600000000001, 0.295968, 0.01009304870291239, 0.01488643315542961, 3, 1, 1.913051666666667, 5.739155, 0.3237647697820488, 0.3383144003623209, 0.1887340252266667, 1.69860622704, 0.003267773589697437, 0.005036294706512939),
    @("FAPs", "Reln", "Lrp8", "sCs", 3, 1, 0.09865600000000001, 0.295968, 0.01009304870291239, 0.01488643315542961, 2, 1, 0.762342, 1.524684, 0.1290187225080277, 0.0898777874446717, 0.07520961235200001, 0.451257674112, 0.001302192249861063, 0.001337959674953016),
    @("FAPs", "Reln", "Lrp8", "ECs", 3, 1, 0.2337746666666667, 0.7013240000000001, 0.02391642775070728, 0.03527480283780177, 2, 0.6666666666666666, 0.7407323333333333, 2.222197, 0.125361503586392, 0.1309951108729331, 0.1731644543142222, 1.558480088828, 0.002998199343243975, 0.004620826708758698),
    @("FAPs", "Reln", "Lrp8", "FAPs", 3, 1, 0.2337746666666667, 0.7013240000000001, 0.02391642775070728, 0.03527480283780177, 3, 1, 1.080588666666667, 3.241766, 0.1828787726899297, 0.1910971424199137, 0.2526142553537778, 2.273528298184, 0.004373806954176722, 0.006740914021729778),
    @("FAPs", "Reln", "Lrp8", "M1", 3, 1, 0.2337746666666667, 0.7013240000000001, 0.02391642775070728, 0.03527480283780177, 3, 1, 1.412055666666667, 4.236167, 0.2389762314336017, 0.2497155589001607, 0.3301028427897778, 2.970925585108, 0.005715457773218037, 0.008808667105734642),
    @("FAPs", "Reln", "Lrp8", "M2", 3, 1, 0.2337746666666667, 0.7013240000000001, 0.02391642775070728, 0.03527480283780177, 3, 1, 1.913051666666667, 5.739155, 0.3237647697820488, 0.3383144003623209, 0.4472230156911112, 4.025007141220001, 0.007743296724716746, 0.01193397376997),
    @("FAPs", "Reln", "Lrp8", "sCs", 3, 1, 0.2337746666666667, 0.7013240000000001, 0.02391642775070728, 0.03527480283780177, 2, 1, 0.762342, 1.524684, 0.1290187225080277, 0.0898777874446717, 0.178216246936, 1.069297481616, 0.003085666955351796, 0.003170421231608649),
    @("sCs", "Reln", "Lrp8", "ECs", 2, 1, 9.4422175, 18.884435, 0.9659905235463803, 0.9498387640067686, 2, 0.6666666666666666, 0.7407323333333333, 2.222197, 0.125361503586392, 0.1309951108729331, 6.994155800615833, 41.964934803695, 0.1210980244819802, 0.1244242342024764),
    @("sCs", "Reln", "Lrp8", "FAPs", 2, 1, 9.4422175, 18.884435, 0.9659905235463803, 0.9498387640067686, 3, 1, 1.080588666666667, 3.241766, 0.1828787726899297, 0.1910971424199137, 10.20315321870167, 61.21891931221001, 0.1766591613762646, 0.1815114735613562),
    @("sCs", "Reln", "Lrp8", "M1", 2, 1, 9.4422175, 18.884435, 0.9659905235463803, 0.9498387640067686, 3, 1, 1.412055666666667, 4.236167, 0.2389762314336017, 0.2497155589001607, 13.33293672677417, 79.997620360645, 0.2308487749176858, 0.2371895178189881),
    @("sCs", "Reln", "Lrp8", "M2", 2, 1, 9.4422175, 18.884435, 0.9659905235463803, 0.9498387640067686, 3, 1, 1.913051666666667, 5.739155, 0.3237647697820488, 0.3383144003623209, 18.06344992540417, 108.380699552425, 0.3127536994676346, 0.3213441318858379),
    @("sCs", "Reln", "Lrp8", "sCs", 2, 1, 9.4422175, 18.884435, 0.9659905235463803, 0.9498387640067686, 2, 1, 0.762342, 1.524684, 0.1290187225080277, 0.0898777874446717, 7.198198973385, 28.79279589354, 0.1246308633028149, 0.08536940653811004)

)

$rowIndex = 2
foreach ($row in $data) {
    $colIndex = 1
    foreach ($val in $row) {
        $ws.Cells.Item($rowIndex, $colIndex).Value = $val
        $colIndex++
    }
    $rowIndex++
}
